$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: mark Diego (7), Natalia (8) and Jose (9) as "Usado" ---
$ws.Range("K7").Value = "Usado"
$ws.Range("K8").Value = "Usado"
$ws.Range("K9").Value = "Usado"

# --- New employee records (rows 14-22) ---
# Field order mirrors the original creation logic: birth date (G), license
# expiration date (F) and first name (A) are written per-row across the whole
# batch before the remaining, already-known values are filled in.

$birthDates = @("2001-04-24","2001-05-24","2001-06-24","2001-04-25","2001-05-25","2001-06-25","2001-04-26","2001-05-26","2001-06-26")
for ($i = 0; $i -lt 9; $i++) {
    $row = 14 + $i
    $ws.Range("G$row").Value = "'" + $birthDates[$i]
}

$licenseExp = @("2019-02-24","2019-04-24","2019-05-24","2019-06-24","2019-02-25","2019-04-25","2019-05-25","2019-06-25","2019-02-26")
for ($i = 0; $i -lt 9; $i++) {
    $row = 14 + $i
    $ws.Range("F$row").Value = "'" + $licenseExp[$i]
}

$firstNames = @("Lorena","Jimmy","Patricia","Melissa","Joaquin","Pilar","Lourdes","Joseph","Melina")
for ($i = 0; $i -lt 9; $i++) {
    $row = 14 + $i
    $ws.Range("A$row").Value = $firstNames[$i]
}

$middleNames = @("Daniela","Carlos","Andrea","Daniela","Carlos","Andrea","Daniela","Carlos","Andrea")
$lastNames   = @("Rios Duque","Diaz Perez","Martinez Castro","Rios Duque","Diaz Perez","Martinez Castro","Rios Duque","Diaz Perez","Martinez Castro")
$cedulas     = @(105369886,105369887,105369888,105369889,105369890,105369891,105369892,105369893,105369894)
$licencias   = @(555563,555564,555565,555566,555567,555568,555569,555570,555571)
$nationalities = @("Colombian","Angolan","Bolivian","Colombian","Angolan","Bolivian","Colombian","Angolan","Bolivian")
$maritalStatus = @("Single","Married","Other","Single","Married","Other","Single","Married","Other")
$genders       = @("Mujer","Hombre","Mujer","Mujer","Hombre","Mujer","Mujer","Hombre","Mujer")
$estados       = @("Disponible","Disponible","Disponible","Disponible","Disponible","Disponible","Disponible","Disponible","Disponible")

for ($i = 0; $i -lt 9; $i++) {
    $row = 14 + $i
    $ws.Range("B$row").Value = $middleNames[$i]
    $ws.Range("C$row").Value = $lastNames[$i]
    $ws.Range("D$row").Value = $cedulas[$i]
    $ws.Range("E$row").Value = $licencias[$i]
    $ws.Range("H$row").Value = "'" + $nationalities[$i]
    $ws.Range("I$row").Value = "'" + $maritalStatus[$i]
    $ws.Range("J$row").Value = "'" + $genders[$i]
    $ws.Range("K$row").Value = $estados[$i]
}

# Column K (Estado) now holds the longest values in the sheet; widen it to fit.
$ws.Columns("K").EntireColumn.AutoFit() | Out-Null

# Keep the same "next empty row selected" behaviour the workbook had before.
$ws.Range("A23").Select() | Out-Null
